$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 27; this shifts the existing rows 27..58
# down to 28..59 and copies the formatting of row 27 (e.g. the date
# style on column D) into the freshly inserted row, matching Excel's
# native "insert row" behaviour.
$ws.Rows.Item(27).EntireRow.Insert()

# Populate the newly inserted row 27 with the new record.
$ws.Range("A27").Value = 8
$ws.Range("B27").Value = "Terminal La Palmera de La Serena"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44554
$ws.Range("E27").Value = 4
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100109
$ws.Range("H27").Value = "Uva"
$ws.Range("I27").Value = 100109001
$ws.Range("J27").Value = "Uva"
$ws.Range("K27").Value = "Flame Seedless"
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 600
$ws.Range("N27").Value = 7500
$ws.Range("O27").Value = 8000
$ws.Range("P27").Value = 7750
$ws.Range("Q27").Value = "$/bandeja 10 kilos"
$ws.Range("R27").Value = "Provincia de Limarí"
$ws.Range("S27").Value = 775
$ws.Range("T27").Value = 10
